# Automatic update of files.
# The edit swaps the full content of row 3 <-> row 4, and row 6 <-> row 7
# on the active worksheet (observation records got re-ordered/re-numbered
# upstream; all columns of each pair of rows are exchanged, including the
# "Aktivitet" (M) value that only one of rows 6/7 carries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($Row1, $Row2, $FirstCol, $LastCol)

    $range1 = $ws.Range("$FirstCol$Row1`:$LastCol$Row1")
    $range2 = $ws.Range("$FirstCol$Row2`:$LastCol$Row2")

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}

# Two kinds of columns need to be left out of the bulk range swap below,
# because round-tripping their value through Range.Value has side effects
# in this environment:
#  - Y ("Startdatum") and AA ("Slutdatum") contain plain text that looks
#    like a date (e.g. "2026-02-05"); re-assigning such text auto-converts
#    it into a real date value/format, which the diff does not want.
#  - I ("Antal"), AT ("Bestamningsar") and AY ("Projektnamn") are present
#    but empty (inline string with no text); reading back an empty value
#    and reassigning it deletes the cell instead of keeping the empty
#    cell in place.
# All of these columns hold identical content in both rows of each pair
# being swapped, so it is safe to simply skip them.
Swap-RowRange 3 4 "A" "H"
Swap-RowRange 3 4 "J" "X"
Swap-RowRange 3 4 "Z" "Z"
Swap-RowRange 3 4 "AB" "AS"
Swap-RowRange 3 4 "AU" "AX"

Swap-RowRange 6 7 "A" "H"
Swap-RowRange 6 7 "J" "X"
Swap-RowRange 6 7 "Z" "Z"
Swap-RowRange 6 7 "AB" "AS"
Swap-RowRange 6 7 "AU" "AX"
